# Re-sort the Madrid neighbourhood average-rating table in descending order
# of averageRating, add the two new "100" neighbourhoods (Valdemarin, Horcajo)
# at the top, and rename the header from "neighbourhood_cleansed" to "neighbourhood".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value2 = "neighbourhood"
$ws.Range("B1").Value2 = "averageRating"

# Data rows (2 .. 128), sorted by averageRating descending
$data = New-Object 'object[,]' 127,2
$data[0,0] = "Valdemarín"
$data[0,1] = 100.0
$data[1,0] = "Horcajo"
$data[1,1] = 100.0
$data[2,0] = "Hellín"
$data[2,1] = 98.58333333333333
$data[3,0] = "Amposta"
$data[3,1] = 98.22222222222223
$data[4,0] = "Corralejos"
$data[4,1] = 98.0
$data[5,0] = "Ambroz"
$data[5,1] = 97.16666666666667
$data[6,0] = "Peñagrande"
$data[6,1] = 97.02702702702703
$data[7,0] = "Campamento"
$data[7,1] = 96.94117647058823
$data[8,0] = "Estrella"
$data[8,1] = 96.78571428571429
$data[9,0] = "Mirasierra"
$data[9,1] = 96.65384615384616
$data[10,0] = "Media Legua"
$data[10,1] = 96.3076923076923
$data[11,0] = "La Paz"
$data[11,1] = 96.1875
$data[12,0] = "Santa Eugenia"
$data[12,1] = 96.16666666666667
$data[13,0] = "El Goloso"
$data[13,1] = 96.11111111111111
$data[14,0] = "Buenavista"
$data[14,1] = 96.02702702702703
$data[15,0] = "Quintana"
$data[15,1] = 95.90789473684211
$data[16,0] = "Arcos"
$data[16,1] = 95.85185185185185
$data[17,0] = "Timón"
$data[17,1] = 95.62162162162163
$data[18,0] = "Butarque"
$data[18,1] = 95.45454545454545
$data[19,0] = "San Pascual"
$data[19,1] = 95.42307692307692
$data[20,0] = "San Juan Bautista"
$data[20,1] = 95.42105263157895
$data[21,0] = "Alameda de Osuna"
$data[21,1] = 95.3103448275862
$data[22,0] = "El Viso"
$data[22,1] = 95.03571428571429
$data[23,0] = "Colina"
$data[23,1] = 95.0
$data[24,0] = "Palomeras Sureste"
$data[24,1] = 95.0
$data[25,0] = "Nueva España"
$data[25,1] = 94.95555555555555
$data[26,0] = "Salvador"
$data[26,1] = 94.9375
$data[27,0] = "Los Angeles"
$data[27,1] = 94.88235294117646
$data[28,0] = "Aravaca"
$data[28,1] = 94.84848484848484
$data[29,0] = "Marroquina"
$data[29,1] = 94.84615384615384
$data[30,0] = "Valverde"
$data[30,1] = 94.84482758620689
$data[31,0] = "Pilar"
$data[31,1] = 94.725
$data[32,0] = "Castillejos"
$data[32,1] = 94.6063829787234
$data[33,0] = "Simancas"
$data[33,1] = 94.59649122807018
$data[34,0] = "El Pardo"
$data[34,1] = 94.5
$data[35,0] = "Ciudad Universitaria"
$data[35,1] = 94.48387096774194
$data[36,0] = "Argüelles"
$data[36,1] = 94.42233009708738
$data[37,0] = "Costillares"
$data[37,1] = 94.36
$data[38,0] = "Zofío"
$data[38,1] = 94.28571428571429
$data[39,0] = "Concepción"
$data[39,1] = 94.28070175438596
$data[40,0] = "Entrevías"
$data[40,1] = 94.21875
$data[41,0] = "Casa de Campo"
$data[41,1] = 94.2063492063492
$data[42,0] = "Piovera"
$data[42,1] = 94.14285714285714
$data[43,0] = "Almenara"
$data[43,1] = 93.93846153846154
$data[44,0] = "Niño Jesús"
$data[44,1] = 93.84848484848484
$data[45,0] = "Ventas"
$data[45,1] = 93.83185840707965
$data[46,0] = "Pavones"
$data[46,1] = 93.8
$data[47,0] = "Canillas"
$data[47,1] = 93.77464788732394
$data[48,0] = "Lucero"
$data[48,1] = 93.69642857142857
$data[49,0] = "Chopera"
$data[49,1] = 93.64285714285714
$data[50,0] = "Atocha"
$data[50,1] = 93.63636363636364
$data[51,0] = "Cármenes"
$data[51,1] = 93.59375
$data[52,0] = "Ciudad Jardín"
$data[52,1] = 93.58064516129032
$data[53,0] = "Ibiza"
$data[53,1] = 93.57352941176471
$data[54,0] = "Pinar del Rey"
$data[54,1] = 93.50877192982456
$data[55,0] = "Casco Histórico de Vicálvaro"
$data[55,1] = 93.5
$data[56,0] = "Jerónimos"
$data[56,1] = 93.48235294117647
$data[57,0] = "Orcasitas"
$data[57,1] = 93.33333333333333
$data[58,0] = "Palomas"
$data[58,1] = 93.3
$data[59,0] = "Puerta del Angel"
$data[59,1] = 93.28834355828221
$data[60,0] = "Delicias"
$data[60,1] = 93.23870967741935
$data[61,0] = "Pacífico"
$data[61,1] = 93.20886075949367
$data[62,0] = "Prosperidad"
$data[62,1] = 93.18691588785046
$data[63,0] = "Valdefuentes"
$data[63,1] = 93.17647058823529
$data[64,0] = "Goya"
$data[64,1] = 93.15228426395939
$data[65,0] = "Numancia"
$data[65,1] = 93.1484375
$data[66,0] = "Almendrales"
$data[66,1] = 93.13953488372093
$data[67,0] = "Legazpi"
$data[67,1] = 93.10526315789474
$data[68,0] = "Aluche"
$data[68,1] = 93.06896551724138
$data[69,0] = "Almagro"
$data[69,1] = 92.95495495495496
$data[70,0] = "Rios Rosas"
$data[70,1] = 92.93162393162393
$data[71,0] = "Moscardó"
$data[71,1] = 92.7090909090909
$data[72,0] = "Cuatro Caminos"
$data[72,1] = 92.70289855072464
$data[73,0] = "Hispanoamérica"
$data[73,1] = 92.68852459016394
$data[74,0] = "San Fermín"
$data[74,1] = 92.68
$data[75,0] = "Bellas Vistas"
$data[75,1] = 92.62745098039215
$data[76,0] = "Pueblo Nuevo"
$data[76,1] = 92.61111111111111
$data[77,0] = "Aguilas"
$data[77,1] = 92.56521739130434
$data[78,0] = "Arapiles"
$data[78,1] = 92.53543307086615
$data[79,0] = "Valdeacederas"
$data[79,1] = 92.5223880597015
$data[80,0] = "Imperial"
$data[80,1] = 92.47916666666667
$data[81,0] = "Canillejas"
$data[81,1] = 92.46428571428571
$data[82,0] = "Comillas"
$data[82,1] = 92.42105263157895
$data[83,0] = "San Cristobal"
$data[83,1] = 92.41666666666667
$data[84,0] = "Palomeras Bajas"
$data[84,1] = 92.39583333333333
$data[85,0] = "Casco Histórico de Barajas"
$data[85,1] = 92.38709677419355
$data[86,0] = "Vinateros"
$data[86,1] = 92.33333333333333
$data[87,0] = "Embajadores"
$data[87,1] = 92.28368017524645
$data[88,0] = "Recoletos"
$data[88,1] = 92.2741935483871
$data[89,0] = "Casco Histórico de Vallecas"
$data[89,1] = 92.24137931034483
$data[90,0] = "Pradolongo"
$data[90,1] = 92.22727272727273
$data[91,0] = "Guindalera"
$data[91,1] = 92.16867469879519
$data[92,0] = "Cortes"
$data[92,1] = 92.16347569955818
$data[93,0] = "Justicia"
$data[93,1] = 92.16216216216216
$data[94,0] = "Acacias"
$data[94,1] = 92.15748031496064
$data[95,0] = "Palacio"
$data[95,1] = 92.14148471615721
$data[96,0] = "Castilla"
$data[96,1] = 92.125
$data[97,0] = "Gaztambide"
$data[97,1] = 92.11864406779661
$data[98,0] = "Palos de Moguer"
$data[98,1] = 92.0954356846473
$data[99,0] = "Vallehermoso"
$data[99,1] = 92.03448275862068
$data[100,0] = "Universidad"
$data[100,1] = 91.82299012693935
$data[101,0] = "Sol"
$data[101,1] = 91.81263858093126
$data[102,0] = "Lista"
$data[102,1] = 91.81132075471699
$data[103,0] = "Puerta Bonita"
$data[103,1] = 91.80392156862744
$data[104,0] = "Berruguete"
$data[104,1] = 91.79452054794521
$data[105,0] = "Castellana"
$data[105,1] = 91.76068376068376
$data[106,0] = "Valdezarza"
$data[106,1] = 91.5
$data[107,0] = "San Isidro"
$data[107,1] = 91.33333333333333
$data[108,0] = "Portazgo"
$data[108,1] = 91.32
$data[109,0] = "Rosas"
$data[109,1] = 91.25
$data[110,0] = "Opañel"
$data[110,1] = 91.24657534246575
$data[111,0] = "Abrantes"
$data[111,1] = 91.2
$data[112,0] = "Trafalgar"
$data[112,1] = 90.62008733624454
$data[113,0] = "Vista Alegre"
$data[113,1] = 90.34285714285714
$data[114,0] = "San Diego"
$data[114,1] = 90.23684210526316
$data[115,0] = "Aeropuerto"
$data[115,1] = 90.16666666666667
$data[116,0] = "Fontarrón"
$data[116,1] = 89.47368421052632
$data[117,0] = "Rejas"
$data[117,1] = 89.37735849056604
$data[118,0] = "Fuente del Berro"
$data[118,1] = 89.29885057471265
$data[119,0] = "San Andrés"
$data[119,1] = 88.56
$data[120,0] = "Adelfas"
$data[120,1] = 88.28813559322033
$data[121,0] = "Orcasur"
$data[121,1] = 88.05
$data[122,0] = "Fuentelareina"
$data[122,1] = 88.0
$data[123,0] = "Los Rosales"
$data[123,1] = 87.4
$data[124,0] = "Apostol Santiago"
$data[124,1] = 87.07692307692308
$data[125,0] = "Cuatro Vientos"
$data[125,1] = 85.15384615384616
$data[126,0] = "El Plantío"
$data[126,1] = 83.0

$ws.Range("A2:B128").Value2 = $data
